# Apply "notes on test data" edits to the taxonomy LCA/RESULT table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 41 (RESULT row for the sv17559 group): fill in the NA taxonomy levels
# plus a note in D41 explaining the NA convention.
$ws.Range("D41").Value = "NA (Eukaryota if no NA as name)"
$ws.Range("E41").Value = "NA"
$ws.Range("F41").Value = "NA"
$ws.Range("G41").Value = "NA"
$ws.Range("H41").Value = "NA"
$ws.Range("I41").Value = "NA"
$ws.Range("J41").Value = "NA"

# --- Row 46 (RESULT row for the sv19589 group): fill in kingdom + NA levels.
$ws.Range("D46").Value = "Eukaryota"
$ws.Range("E46").Value = "NA"
$ws.Range("F46").Value = "NA"
$ws.Range("G46").Value = "NA"
$ws.Range("H46").Value = "NA"

# --- Notes column (M) on selected "Bayes" rows.
$ws.Range("M48").Value = "ties -> grab the one that user defined"
$ws.Range("M49").Value = "make sure the table name lines up with the table name"

# --- Row 51 (RESULT row for the sv17897 group): fill in taxonomy through class
# plus remaining NA levels.
$ws.Range("D51").Value = "Eukaryota"
$ws.Range("E51").Value = "Alveolata"
$ws.Range("F51").Value = "Apicomplexa"
$ws.Range("G51").Value = "Coccidiomorphea"
$ws.Range("H51").Value = "NA"
$ws.Range("I51").Value = "NA"
$ws.Range("J51").Value = "NA"

# --- Row 61 (RESULT row for the sv104 group): fill in kingdom + NA levels.
$ws.Range("D61").Value = "Eukaryota"
$ws.Range("E61").Value = "NA"
$ws.Range("F61").Value = "NA"
$ws.Range("G61").Value = "NA"
$ws.Range("H61").Value = "NA"
$ws.Range("I61").Value = "NA"

# --- More notes column (M) on a "Bayes" row.
$ws.Range("M63").Value = "automatic tie breaker -> NA"

# --- View state: zoom in and scroll so row 51 is at the top, select J61.
$ws.Range("J61").Select()
$excel.ActiveWindow.Zoom = 138
$excel.ActiveWindow.ScrollRow = 51
$excel.ActiveWindow.ScrollColumn = 2
